# Apply NOAA-updated temperature (column I) and recomputed ASHP COP values
# (columns N/O) for the affected facility rows in the 325110 longform
# electrification options workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> updated values. Only column I (average_county_temperature)
# changes for every affected row; columns N (worst_ashp_cop) and
# O (best_ashp_cop) only change where they already exist on that row.
$updates = @(
    @{ Row = 2;  I = 12.67039049919483;  N = 1.03592820319873;  O = 1.067704754529766 },
    @{ Row = 3;  I = 12.67039049919483 },
    @{ Row = 4;  I = 1.791666666666668 },
    @{ Row = 5;  I = 1.791666666666668; N = 0.9939102066179896; O = 1.022720671292561 },
    @{ Row = 6;  I = 1.791666666666668 },
    @{ Row = 32; I = 20.68981481481483 },
    @{ Row = 33; I = 19.65277777777778; N = 1.064821081830791;  O = 1.098722912453048 },
    @{ Row = 34; I = 19.65277777777778 },
    @{ Row = 35; I = 20.68981481481483; N = 1.069250338898071;  O = 1.103484165522044 },
    @{ Row = 36; I = 20.68981481481483 },
    @{ Row = 49; I = 13.75752314814816 },
    @{ Row = 51; I = 20.22222222222222 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("I$r").Value = $u.I
    if ($u.ContainsKey("N")) {
        $ws.Range("N$r").Value = $u.N
    }
    if ($u.ContainsKey("O")) {
        $ws.Range("O$r").Value = $u.O
    }
}
